# Adds a new "Italy" market card sheet (copied from "Slovakia"), removes the
# "FBI800" color-code row that Italy doesn't have, fills in the Italy-specific
# values, and makes "Italy" the active/selected sheet (matching what Excel
# records when a user duplicates a sheet and finishes editing it).

$wb = $excel.ActiveWorkbook

# The previously-active sheet ("Slovakia") loses its special "tabSelected"
# view state once a new sheet becomes active; Excel also drops its old
# cell-range selection in favor of a "select all" state recorded for the
# sheet that was active right before the switch.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Activate()
$slovakia.Cells.Select()

# Duplicate the Slovakia sheet (keeps all formatting / merged cells / styles)
# and place the copy right after it, then rename it to "Italy".
$slovakia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Slovakia's card has an extra "FBI800" color-code row that Italy's card
# doesn't include, so remove it (rows below shift up).
$italy.Rows("8").Delete()

# Fill in the Italy-specific market name and NGC reference.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2165"

# Leave the new sheet's selection/view state matching a freshly-edited sheet.
$italy.Range("F21").Select()
